$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "RF04. O sistema deve disponibilizar perguntas para entender o nível do usuário.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF04. O sistema deve disponibilizar um formulário a ser respondido pelo usuário ao logar na plataforma pela primeira vez.",
    2)

$d.Content.Find.Execute(
    "RF05. O sistema deve obter do usuário informações básicas sobre os seus conhecimentos atuais.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF05. O sistema deve obter do usuário informações básicas sobre seus conhecimentos atuais e interesses em programação.",
    2)

$d.Content.Find.Execute(
    "RF10. Os materiais de estudo devem ser rotulados no seu tipo. ex. : vídeo, artigo, livro.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF10. Os materiais de estudo devem ser rotulados no seu tipo. ex. : texto, vídeo, artigo, livro.",
    2)

$d.Content.Find.Execute(
    "RF11. O usuário deve possuir uma forma de marcar o subtópico sugerido ao estudo como completo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF11. O usuário deve possuir uma forma de marcar o subtópico sugerido ao estudo como concluído.",
    2)

$d.Content.Find.Execute(
    "RF12. Subtópicos completos devem recompensar o usuário com pontos de experiência.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF12. Subtópicos concluídos devem recompensar o usuário com pontos de experiência.",
    2)

$d.Content.Find.Execute(
    "Caso o usuário troque a jornada, sua jornada antiga será mantida salva na sua conta. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Caso o usuário opte por uma outra jornada, sua jornada sugerida pelo sistema será mantida na sua conta. ",
    2)

$d.Content.Find.Execute(
    "RF22.  Alterar informações na tabela conforme os usuários vão subindo de pontuação.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF22.  Alterar informações no ranking conforme os usuários vão aumentando suas pontuações.",
    2)

$d.Content.Find.Execute(
    "RF24. Apresentar um botão de salvar cada anotação.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF24. Apresentar um botão para salvar cada anotação.",
    2)

$d.Content.Find.Execute(
    "RF25. Criação de um botão de voltar para o usuário navegar entre as páginas ‘jornadas’.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF25. Criação de um botão para voltar durante a navegação do usuário entre as páginas ‘jornadas’.",
    2)
